$d = $word.ActiveDocument

# Vertical-tab char used by Word for manual line breaks (<w:br/>) inside a run.
$vt = [char]11

# The run of text (with its trailing line break) that needs to move from the
# start of the "Requisitos" bullet list to the end of that list.
$moved = "LOB1019 -  F" + [char]0xED + "sica II  (Requisito fraco)" + $vt

# Locate and remove it from its current position (right after "Requisitos").
$find = $d.Content
$found = $find.Find.Execute($moved, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'LOB1019' requirement line to move"
}
$find.Delete()

# Locate the end of the "LOB1004 -  Cálculo II  (Requisito fraco)" line
# (including its trailing break) and insert the moved text right after it.
$anchor = "LOB1004 -  C" + [char]0xE1 + "lculo II  (Requisito fraco)" + $vt
$target = $d.Content
$found2 = $target.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'LOB1004' requirement line to anchor on"
}
$target.Collapse(0)
$target.InsertAfter($moved)
